# Added updateCard and hooked it up to GUI
#
# TestData.xlsx now tracks when a card was last reviewed: the CARDS sheet
# gets a new LAST_REVIEW column (storing millisecond timestamps - here just
# stubbed with the literal "MILLIS" placeholder for every existing row),
# and the workbook opens on the CARDS sheet instead of MOVES.

$wb = $excel.ActiveWorkbook

$cards = $wb.Worksheets.Item("CARDS")

# New column F: header + one value per existing data row (rows 2-11)
$cards.Cells.Item(1, 6).Value = "LAST_REVIEW"
for ($row = 2; $row -le 11; $row++) {
    $cards.Cells.Item($row, 6).Value = "MILLIS"
}

# Widen columns C:E to fit their content now that the sheet has changed
$cards.Columns.Item(3).ColumnWidth = 18.5
$cards.Columns.Item(4).ColumnWidth = 16.6
$cards.Columns.Item(5).ColumnWidth = 14.6

# GUI now opens on CARDS (previously MOVES was the active/selected tab)
$cards.Activate() | Out-Null
$cards.Range("B3").Select() | Out-Null
